$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Longe  pour Porte-badge"
$ws.Range("H3").Value = "Longe porte-badge"

$ws.Range("G3").Select()
